$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to be interpreted as text so Excel does not
    # auto-convert numeric-looking strings (e.g. "1.00" -> 1, "5.60" -> 5.6)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.938.12"
$ws.Range("E2").Value = "  +4.82%  "

Set-TextValue $ws.Range("D3") "3.350.28"
$ws.Range("E3").Value = "  +4.91%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "556.82"
$ws.Range("E5").Value = "  +3.54%  "

Set-TextValue $ws.Range("D6") "152.71"
$ws.Range("E6").Value = "  +5.47%  "

$ws.Range("E7").Value = "  -0.02%  "

Set-TextValue $ws.Range("D8") "0.531"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("E9").Value = "  +2.40%  "

Set-TextValue $ws.Range("D10") "0.119"
$ws.Range("E10").Value = "  +4.33%  "

Set-TextValue $ws.Range("D11") "0.437"
$ws.Range("E11").Value = "  +1.66%  "

Set-TextValue $ws.Range("D12") "3.928.64"
$ws.Range("E12").Value = "  +4.88%  "

$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("E14").Value = "  +4.24%  "

Set-TextValue $ws.Range("D15") "26.95"
$ws.Range("E15").Value = "  +2.95%  "

Set-TextValue $ws.Range("D16") "62.953.08"
$ws.Range("E16").Value = "  +4.76%  "

Set-TextValue $ws.Range("D17") "3.346.16"
$ws.Range("E17").Value = "  +4.58%  "

Set-TextValue $ws.Range("D18") "6.48"
$ws.Range("E18").Value = "  +4.37%  "

$ws.Range("E19").Value = "  +5.02%  "

$ws.Range("E20").Value = "  +0.69%  "

Set-TextValue $ws.Range("D21") "389.11"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D22") "0.542"
$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.30%  "

Set-TextValue $ws.Range("D24") "70.69"
$ws.Range("E24").Value = "  +0.56%  "

Set-TextValue $ws.Range("D25") "0.181"
$ws.Range("E25").Value = "  +5.74%  "

Set-TextValue $ws.Range("D26") "8.83"
$ws.Range("E26").Value = "  +0.57%  "

Set-TextValue $ws.Range("D27") "0.0₃0974"
$ws.Range("E27").Value = "  +8.44%  "

Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.15%  "

Set-TextValue $ws.Range("D29") "1.99"
$ws.Range("E29").Value = "  +4.31%  "

Set-TextValue $ws.Range("D30") "6.42"
$ws.Range("E30").Value = "  +4.19%  "

Set-TextValue $ws.Range("D31") "23.04"
$ws.Range("E31").Value = "  +3.07%  "

Set-TextValue $ws.Range("D32") "5.60"
$ws.Range("E32").Value = "  +4.10%  "

Set-TextValue $ws.Range("D33") "1.31"
$ws.Range("E33").Value = "  +7.96%  "

Set-TextValue $ws.Range("D34") "6.71"
$ws.Range("E34").Value = "  +2.90%  "

$ws.Range("E35").Value = "  +9.88%  "

Set-TextValue $ws.Range("D36") "159.54"
$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("E37").Value = "  +12.02%  "

Set-TextValue $ws.Range("D38") "27.37"
$ws.Range("E38").Value = "  +6.98%  "

Set-TextValue $ws.Range("D39") "0.0748"
$ws.Range("E39").Value = "  +5.06%  "

Set-TextValue $ws.Range("D40") "2.846.02"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("E41").Value = "  +8.68%  "

Set-TextValue $ws.Range("D42") "4.33"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D43") "40.71"
$ws.Range("E43").Value = "  +2.39%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D44") "0.746"
$ws.Range("E44").Value = "  +2.43%  "

$ws.Range("E45").Value = "  +3.30%  "

Set-TextValue $ws.Range("D46") "3.394.50"
$ws.Range("E46").Value = "  +4.91%  "

Set-TextValue $ws.Range("D47") "22.05"
$ws.Range("E47").Value = "  +7.78%  "

$ws.Range("E48").Value = "  +2.75%  "

Set-TextValue $ws.Range("D49") "6.30"
$ws.Range("E49").Value = "  +1.79%  "

Set-TextValue $ws.Range("D50") "0.807"
$ws.Range("E50").Value = "  +1.41%  "

Set-TextValue $ws.Range("D51") "282.98"
$ws.Range("E51").Value = "  +7.79%  "
